# Progress Log.xlsx edit script
# - Adds 4 new "Week N" labels (Week 3-6) and 6 new progress description strings
# - Adds rows 175-180 to Sheet1 continuing the Summer 2015 weekly log
# - Updates selection to the new last entry cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of writes controls the order new entries are appended to the shared
# string table, so the new "Week N" labels are written first (column A for
# rows 177-180), followed by column A for rows 175-176 (these reuse the
# existing "Week 1"/"Week 2" strings), and finally the column B descriptions
# in the order needed to reproduce the target shared string table layout.

$ws.Range("A177").Value = "Week 3"
$ws.Range("A178").Value = "Week 4"
$ws.Range("A179").Value = "Week 5"
$ws.Range("A180").Value = "Week 6"

$ws.Range("A175").Value = "Week 1"
$ws.Range("A176").Value = "Week 2"

$ws.Range("B175").Value = "Put app on VC; importing contacts from AB; friend suggestions; contactsVC alphabetized"
$ws.Range("B177").Value = "Scrollable button label; push notif/AB permissions flow; interv. Calvin; posted to FB groups"
$ws.Range("B178").Value = "Changed searchVC UI; new launch screens; general debugging/testing; interv. Sean"
$ws.Range("B179").Value = "Master links; screenshots; offer/contract for Sean; submitted archive; created app video"
$ws.Range("B176").Value = "ContactsVC - sending texts to non-users, Recents section, icons and UI; LinkedIn inMail"
$ws.Range("B180").Value = "Refined video; helped Sean; email pitch; Dropbox press kit; emails to 19 journalists (+7)"

# Scroll the view down to the new rows and select the next empty cell below
# the newly added data, matching where the author's cursor ended up.
$excel.Goto($ws.Range("A149"), $true)
$ws.Range("B181").Select()
